# PIN-847: improve handling unresolvable duplicities
#
# Adds a "National ID" (ID Type / ID Number) entry for the three
# beneficiary rows on the "Worksheet" sheet, and updates the sheet's
# stored selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 (Maya Dixon) - National ID 132
$ws.Cells.Item(7, 37).Value = "National ID"
$ws.Cells.Item(7, 38).Value = 132

# Row 8 (Ashton Harris) - National ID 88
$ws.Cells.Item(8, 37).Value = "National ID"
$ws.Cells.Item(8, 38).Value = 88

# Row 9 (Abigail Craig) - National ID 61561
$ws.Cells.Item(9, 37).Value = "National ID"
$ws.Cells.Item(9, 38).Value = 61561

# Move the active selection to AL10, matching the updated sheet view.
$ws.Range("AL10").Select() | Out-Null
